# Automatic update of files.
# The "Förändrad" (Changed) column (C) for every data row advances by one
# day, from serial date 46074 (2026-02-21) to 46075 (2026-02-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
